# Adjusting the NRG location
# Shift all timestamps in column A (rows 2-97) forward by one day and
# update the solar production values in column B for the rows affected
# by the shifted production curve.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateValues = @{
    2 = 45947.01041666666
    3 = 45947.02083333334
    4 = 45947.03125
    5 = 45947.04166666666
    6 = 45947.05208333334
    7 = 45947.0625
    8 = 45947.07291666666
    9 = 45947.08333333334
    10 = 45947.09375
    11 = 45947.10416666666
    12 = 45947.11458333334
    13 = 45947.125
    14 = 45947.13541666666
    15 = 45947.14583333334
    16 = 45947.15625
    17 = 45947.16666666666
    18 = 45947.17708333334
    19 = 45947.1875
    20 = 45947.19791666666
    21 = 45947.20833333334
    22 = 45947.21875
    23 = 45947.22916666666
    24 = 45947.23958333334
    25 = 45947.25
    26 = 45947.26041666666
    27 = 45947.27083333334
    28 = 45947.28125
    29 = 45947.29166666666
    30 = 45947.30208333334
    31 = 45947.3125
    32 = 45947.32291666666
    33 = 45947.33333333334
    34 = 45947.34375
    35 = 45947.35416666666
    36 = 45947.36458333334
    37 = 45947.375
    38 = 45947.38541666666
    39 = 45947.39583333334
    40 = 45947.40625
    41 = 45947.41666666666
    42 = 45947.42708333334
    43 = 45947.4375
    44 = 45947.44791666666
    45 = 45947.45833333334
    46 = 45947.46875
    47 = 45947.47916666666
    48 = 45947.48958333334
    49 = 45947.5
    50 = 45947.51041666666
    51 = 45947.52083333334
    52 = 45947.53125
    53 = 45947.54166666666
    54 = 45947.55208333334
    55 = 45947.5625
    56 = 45947.57291666666
    57 = 45947.58333333334
    58 = 45947.59375
    59 = 45947.60416666666
    60 = 45947.61458333334
    61 = 45947.625
    62 = 45947.63541666666
    63 = 45947.64583333334
    64 = 45947.65625
    65 = 45947.66666666666
    66 = 45947.67708333334
    67 = 45947.6875
    68 = 45947.69791666666
    69 = 45947.70833333334
    70 = 45947.71875
    71 = 45947.72916666666
    72 = 45947.73958333334
    73 = 45947.75
    74 = 45947.76041666666
    75 = 45947.77083333334
    76 = 45947.78125
    77 = 45947.79166666666
    78 = 45947.80208333334
    79 = 45947.8125
    80 = 45947.82291666666
    81 = 45947.83333333334
    82 = 45947.84375
    83 = 45947.85416666666
    84 = 45947.86458333334
    85 = 45947.875
    86 = 45947.88541666666
    87 = 45947.89583333334
    88 = 45947.90625
    89 = 45947.91666666666
    90 = 45947.92708333334
    91 = 45947.9375
    92 = 45947.94791666666
    93 = 45947.95833333334
    94 = 45947.96875
    95 = 45947.97916666666
    96 = 45947.98958333334
    97 = 45948
}

foreach ($row in $dateValues.Keys) {
    $ws.Cells.Item($row, 1).Value2 = $dateValues[$row]
}

# Updated "Actual Production (MW)" values for rows 29-43 to reflect the
# production curve for the new day.
$newValues = @{
    29 = 11
    30 = 49
    31 = 106
    32 = 185
    33 = 289
    34 = 421
    35 = 541
    36 = 654
    37 = 761
    38 = 828
    39 = 916
    40 = 982
    41 = 1051
    42 = 1124
    43 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
